$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.443.86'
$ws.Range('E2').Value = '  +5.78%  '
$ws.Range('D3').Value = '3.549.74'
$ws.Range('E3').Value = '  +1.99%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '''417.43'
$ws.Range('E5').Value = '  +0.36%  '
$ws.Range('D6').Value = '''129.64'
$ws.Range('E6').Value = '  -0.14%  '
$ws.Range('D7').Value = '''0.652'
$ws.Range('E7').Value = '  +4.13%  '
$ws.Range('D8').Value = '3.539.78'
$ws.Range('E8').Value = '  +1.88%  '
$ws.Range('D9').Value = '''1.00'
$ws.Range('E9').Value = '  +0.01%  '
$ws.Range('D10').Value = '''0.776'
$ws.Range('E10').Value = '  +6.08%  '
$ws.Range('E11').Value = '  +28.81%  '
$ws.Range('D12').Value = '''0.0000356'
$ws.Range('E12').Value = '  +62.61%  '
$ws.Range('D13').Value = '''42.77'
$ws.Range('E13').Value = '  -0.40%  '
$ws.Range('D14').Value = '''10.00'
$ws.Range('E14').Value = '  +4.16%  '
$ws.Range('D15').Value = '4.118.90'
$ws.Range('E15').Value = '  +2.45%  '
$ws.Range('E16').Value = '  -0.27%  '
$ws.Range('D17').Value = '''20.30'
$ws.Range('E17').Value = '  -1.57%  '
$ws.Range('D18').Value = '3.537.56'
$ws.Range('E18').Value = '  +1.60%  '
$ws.Range('E19').Value = '  +4.06%  '
$ws.Range('D20').Value = '''12.48'
$ws.Range('E20').Value = '  -3.21%  '
$ws.Range('D21').Value = '66.348.80'
$ws.Range('E21').Value = '  +5.54%  '
$ws.Range('D22').Value = '''446.48'
$ws.Range('E22').Value = '  -5.42%  '
$ws.Range('D23').Value = '''89.68'
$ws.Range('E23').Value = '  -1.67%  '
$ws.Range('D24').Value = '''3.19'
$ws.Range('E24').Value = '  -3.29%  '
$ws.Range('D25').Value = '''13.10'
$ws.Range('E25').Value = '  -2.47%  '
$ws.Range('D26').Value = '''3.36'
$ws.Range('E26').Value = '  +1.23%  '
$ws.Range('D27').Value = '''9.94'
$ws.Range('E27').Value = '  -6.06%  '
$ws.Range('D28').Value = '''34.36'
$ws.Range('E28').Value = '  +2.41%  '
$ws.Range('E29').Value = '  +0.39%  '
$ws.Range('D31').Value = '''12.40'
$ws.Range('E31').Value = '  +3.22%  '
$ws.Range('E32').Value = '  +3.60%  '
$ws.Range('D33').Value = '''7.24'
$ws.Range('E33').Value = '  -5.44%  '
$ws.Range('E34').Value = '  -4.60%  '
$ws.Range('E35').Value = '  -0.16%  '
$ws.Range('D36').Value = '''38.90'
$ws.Range('E36').Value = '  -4.76%  '
$ws.Range('D37').Value = '0.0₃0810'
$ws.Range('E37').Value = '  +45.44%  '
$ws.Range('D38').Value = '''56.75'
$ws.Range('E38').Value = '  -2.33%  '
$ws.Range('D39').Value = '''0.0495'
$ws.Range('E39').Value = '  +0.68%  '
$ws.Range('D40').Value = '''0.146'
$ws.Range('E40').Value = '  +8.50%  '
$ws.Range('D41').Value = '''0.999'
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('D42').Value = '''2.77'
$ws.Range('E42').Value = '  +2.93%  '
$ws.Range('D43').Value = '''3.00'
$ws.Range('E43').Value = '  -0.99%  '
$ws.Range('D44').Value = '''148.50'
$ws.Range('E44').Value = '  +2.13%  '
$ws.Range('D45').Value = '''4.37'
$ws.Range('E45').Value = '  -0.42%  '
$ws.Range('D46').Value = '''3.23'
$ws.Range('E46').Value = '  -4.32%  '
$ws.Range('D47').Value = '''0.307'
$ws.Range('E47').Value = '  -5.83%  '
$ws.Range('E48').Value = '  -5.53%  '
$ws.Range('D49').Value = '''2.36'
$ws.Range('E49').Value = '  -1.92%  '
$ws.Range('B50').Value = 'ApeXProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D50').Value = '''2.60'
$ws.Range('E50').Value = '  +10.98%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = '''0.143'
$ws.Range('E51').Value = '  +2.27%  '
